# Updated account creation page
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Email" sheet: fix the sample email address + matching hyperlink
# ---------------------------------------------------------------------
$wsEmail = $wb.Worksheets.Item("Email")
$wsEmail.Range("A2").Value = "anilkumar@gmail.com"

$wsEmail.Hyperlinks.Delete()
$wsEmail.Hyperlinks.Add($wsEmail.Range("A2"), "mailto:anilkumar@gmail.com", "", "", "anilkumar@gmail.com")
# Hyperlinks.Add() re-styles the cell with a generic "Hyperlink" style; put the
# original (already-hyperlinked) look back so we don't fork a new cell style.
$fEmail = $wsEmail.Range("A2").Font
$fEmail.Name = "Calibri"
$fEmail.Size = 11
$fEmail.Underline = 2
$fEmail.Color = 16711680

# ---------------------------------------------------------------------
# 2. "AccountCreationData" sheet: refresh the three sample rows
# ---------------------------------------------------------------------
$wsAcc = $wb.Worksheets.Item("AccountCreationData")

# Row 2
$wsAcc.Range("A2").Value = "fsfsfdsf@gmail.com"
$wsAcc.Range("E2").Value = "test@123"
$wsAcc.Range("I2").Value = "meslova"
$wsAcc.Range("J2").Value = "hytech"
$wsAcc.Range("O2").Value = "9556238894"

# Row 3
$wsAcc.Range("A3").Value = "agwv@gmail.com"
$wsAcc.Range("E3").Value = "test@124"
$wsAcc.Range("I3").Value = "tcs"
$wsAcc.Range("J3").Value = "madhpur"
$wsAcc.Range("O3").Value = "7077777607"

# Row 4
$wsAcc.Range("A4").Value = "dfdagwv@gmail.com"
$wsAcc.Range("E4").Value = "test@125"
$wsAcc.Range("I4").Value = "cts"
$wsAcc.Range("J4").Value = "ammerpet"
$wsAcc.Range("O4").Value = "9777997046"

# Rebuild the hyperlinks: A2, E2 (new), A3, A4
$wsAcc.Hyperlinks.Delete()
$wsAcc.Hyperlinks.Add($wsAcc.Range("A2"), "mailto:fsfsfdsf@gmail.com", "", "", "fsfsfdsf@gmail.com")
$wsAcc.Hyperlinks.Add($wsAcc.Range("E2"), "mailto:test@123", "", "", "test@123")
$wsAcc.Hyperlinks.Add($wsAcc.Range("A3"), "mailto:agwv@gmail.com", "", "", "agwv@gmail.com")
$wsAcc.Hyperlinks.Add($wsAcc.Range("A4"), "mailto:dfdagwv@gmail.com", "", "", "dfdagwv@gmail.com")

# Restore original (non-"new hyperlink default") look on the linked cells
foreach ($addr in @("A2", "E2", "A3", "A4")) {
    $f = $wsAcc.Range($addr).Font
    $f.Name = "Calibri"
    $f.Size = 11
    $f.Underline = 2
    $f.Color = 16711680
}

# Row heights on the data rows changed slightly when the sheet was re-saved
$wsAcc.Rows.Item(2).RowHeight = 13.8
$wsAcc.Rows.Item(3).RowHeight = 13.8
$wsAcc.Rows.Item(4).RowHeight = 13.8

# ---------------------------------------------------------------------
# 3. Active sheet / selection bookkeeping
# ---------------------------------------------------------------------
$wsEmail.Activate()
$wsEmail.Range("C9").Select()

$wsAcc.Activate()
$wsAcc.Range("J5").Select()
